# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# -------- Sheet "展览" (sheetId 1) --------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value2 = 4668
$wsExpo.Range("F10").Value2 = 197
$wsExpo.Range("F11").Value2 = 180
$wsExpo.Range("F12").Value2 = 1756
$wsExpo.Range("F14").Value2 = 3933
$wsExpo.Range("F16").Value2 = 268

# -------- Sheet "演出" (sheetId 2) --------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value2 = 46

# -------- Sheet "全部类型" (sheetId 4) --------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value2 = 4668
$wsAll.Range("F9").Value2 = 46
$wsAll.Range("F12").Value2 = 197
$wsAll.Range("F13").Value2 = 180
$wsAll.Range("F16").Value2 = 1757
$wsAll.Range("F18").Value2 = 3933
$wsAll.Range("F20").Value2 = 268
